# Generate Report for Handoff
# Updates the localization-status report: the e78a2ec2 and ead2e901 entries
# move from "Handed back: in sync with en-US" to "Ready for handoff", with
# refreshed handoff timestamps and an explanatory "not the latest" error
# detail for the zh-cn / de-de language sheets.

$wb = $excel.ActiveWorkbook

$notLatest_e78a = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2e63be22ffe7659ee0836895ced8c6a6bdc99f4b/e2e/e78a2ec2-34b6-494b-b3f4-6a1f092c7816.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9f0471b89f708deb43ed9c7fcecabb9dcf14eefb/e2e/e78a2ec2-34b6-494b-b3f4-6a1f092c7816.md."
$notLatest_ead2 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2e63be22ffe7659ee0836895ced8c6a6bdc99f4b/e2e/ead2e901-fe44-4608-9273-b69037b91c7f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9f0471b89f708deb43ed9c7fcecabb9dcf14eefb/e2e/ead2e901-fe44-4608-9273-b69037b91c7f.md."

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-09-02 00:34:58"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-09-02 00:34:58"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("H4").Value = "2016-09-02 00:34:54"
$wsZhCn.Range("P4").Value = $notLatest_e78a
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("H5").Value = "2016-09-02 00:34:54"
$wsZhCn.Range("P5").Value = $notLatest_ead2
# Widen the Error Detail column (P) to fit the new, longer message
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("H4").Value = "2016-09-02 00:34:58"
$wsDeDe.Range("P4").Value = $notLatest_e78a
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("H5").Value = "2016-09-02 00:34:58"
$wsDeDe.Range("P5").Value = $notLatest_ead2
# Widen the Error Detail column (P) to fit the new, longer message
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
